# Append the profit figure for 2025-10-15 as a new row at the bottom of
# the sheet (row 59), matching the pattern of the existing data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 59
$dateText = "10/15/2025"
$profitValue = 11034.92

# Writing a date-shaped string straight into a cell's .Value triggers
# Excel's normal auto-detection, turning it into a serial-number date
# (plus a date NumberFormat) instead of the literal text used by every
# other row in this column. Build the literal text as a formula result
# in a scratch cell first (a formula's string result is never
# reinterpreted as a date), then copy just the *value* into place so the
# destination cell keeps the plain/default formatting of its neighbors.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""$dateText"""
$scratch.Copy()
$ws.Range("A$newRow").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Cells.Item($newRow, 2).Value = $profitValue
